# "Updates to prop never treated and coverage times"
#
# Fills in the "never treated" coverage row (row 2, New Product A / MDA /
# School) on the "Platform Coverage" sheet for every year between 2018 and
# 2040 (previously only alternating years had a value) and refreshes the
# active cell selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")
$ws.Activate() | Out-Null

# Row 2 originally only had coverage = 0.6 on the even-lettered columns
# (H, J, L, N, P, R, T, V, X, Z, AB, AD). Fill in the gaps (I, K, M, O, Q,
# S, U, W, Y, AA, AC) so every year from 2018-2040 has the same 0.6 value.
$coverageCols = @("I", "K", "M", "O", "Q", "S", "U", "W", "Y", "AA", "AC")
foreach ($col in $coverageCols) {
    $ws.Range($col + "2").Value = 0.6
}

# Update the active selection on the sheet.
$ws.Range("AB6").Select() | Out-Null
